$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.137.00'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '2.928.21'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.04'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.73%  '
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.442'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '3.413.16'
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = '61.099.54'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = '2.925.79'
$ws.Range("E18").Value = '  +0.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '432.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("E21").Value = '  -0.48%  '
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.24%  '
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E26").Value = '  -1.56%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -3.22%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("E31").Value = '  +3.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  +3.19%  '
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("E36").Value = '  +0.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.98'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.95%  '
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  -0.66%  '
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '41.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.28%  '
$ws.Range("E42").Value = '  -1.90%  '
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '371.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("D45").Value = '2.697.94'
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.03%  '
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = '  +0.08%  '
